# Applies the "2nd update for user story" edits to the Customer user
# story document. Each table-cell paragraph that changes is rewritten
# via Range.InsertXML so the exact run layout (including w:proofErr
# grammar-check markers that Word's editor would have inserted) matches
# the authored revision.

$d = $word.ActiveDocument

function Set-ParaRuns {
    param(
        $Doc,
        [string]$Needle,
        [string]$InnerXml,
        [bool]$All = $false
    )

    $count = $Doc.Paragraphs.Count
    $changed = 0
    for ($i = 1; $i -le $count; $i++) {
        $p = $Doc.Paragraphs($i)
        if ($p.Range.Text.Contains($Needle)) {
            $pkg = '<?xml version="1.0" standalone="yes"?>' +
                '<?mso-application progid="Word.Document"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p>' + $InnerXml + '</w:p></w:body>' +
                '</w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'
            $p.Range.InsertXML($pkg)
            $changed = $changed + 1
            if (-not $All) {
                return $changed
            }
        }
    }
    return $changed
}

# 1) "To notify marina of arrival " -> split so "marina" is flagged
#    with a grammar-check proofErr pair.
Set-ParaRuns $d "notify marina of" (
    '<w:r><w:t>To</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> notify </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>marina</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> of </w:t></w:r>' +
    '<w:r><w:t>arrival</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
) | Out-Null

# 2) "So Staff will know ..." -> "So" flagged with proofErr pair.
Set-ParaRuns $d "So Staff will know they are needed" (
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>So</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Staff will know they are needed when customers with Docking assistance services arrive and depart.</w:t></w:r>'
) | Out-Null

# 3) "Marina " + "Employee" + ", " + "Cindy Wells" -> single merged run.
#    (Applies to every paragraph with this text; the two that are
#    already a single run are left unchanged in substance.)
Set-ParaRuns $d "Marina Employee, Cindy Wells" (
    '<w:r><w:t>Marina Employee, Cindy Wells</w:t></w:r>'
) $true | Out-Null

# 4) "To set up customer accounts" -> "To manage additional services"
#    with "To manage" flagged by proofErr.
Set-ParaRuns $d "To s" (
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>To manage</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> additional services</w:t></w:r>'
) | Out-Null

# 5) "I can add new customers to the system and offer extra services
#    they might need." -> "I can add or remove extra services for
#    customers that they might need."
Set-ParaRuns $d "add new customers to the system" (
    '<w:r><w:t>I</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> can add or remove </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">extra services </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">for customers that </w:t></w:r>' +
    '<w:r><w:t>they might need.</w:t></w:r>'
) | Out-Null

# 6) "To check slip availability" -> "To find slips by customer name"
Set-ParaRuns $d "To check slip availability" (
    '<w:r><w:t xml:space="preserve">To </w:t></w:r>' +
    '<w:r><w:t>find slips by customer name</w:t></w:r>'
) | Out-Null

# 7) "I can help new customers find a spot for their boat." ->
#    "I can quickly find where the customer's boat is located."
Set-ParaRuns $d "I can help new customers find a spot for their boat." (
    '<w:r><w:t xml:space="preserve">I can </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">quickly find where </w:t></w:r>' +
    '<w:r><w:t>the customer' + [char]0x2019 + 's</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> boat is located</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>'
) | Out-Null

# 8) "To " + "keep up with customer payments" -> single merged run.
Set-ParaRuns $d "keep up with customer payments" (
    '<w:r><w:t>To keep up with customer payments</w:t></w:r>'
) | Out-Null

# 9) "... and take action if needed." -> "take action" flagged by
#    proofErr, rest of the paragraph untouched in substance.
Set-ParaRuns $d "take action if needed" (
    '<w:r><w:t xml:space="preserve">I can </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">check which customers haven' + [char]0x2019 + 't paid yet </w:t></w:r>' +
    '<w:r><w:t>and</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>take action</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> if needed</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>'
) | Out-Null

# 10) "A better way for employees to communicate" ->
#     "To track slip availability"
Set-ParaRuns $d "A better way for employees to communicate" (
    '<w:r><w:t>To track slip availability</w:t></w:r>'
) | Out-Null

# 11) "They can work together more smoothly and provide better
#     customer service." -> "I can offer special prices for slips
#     with lots of vacancies." with "slips" flagged by proofErr.
Set-ParaRuns $d "hey can work together more smoothly" (
    '<w:r><w:t xml:space="preserve">I can offer special prices for </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>slips</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> with lots of vacancies</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>'
) | Out-Null

# 12) "To plan and manage projects" -> "To find customer contact info"
Set-ParaRuns $d "plan and manage projects" (
    '<w:r><w:t>To</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>find</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> customer contact info</w:t></w:r>'
) | Out-Null

# 13) "I can have clear steps before starting and make sure everything
#     runs smoothly without unexpected issues." -> "I can access
#     customer contact details by slip number and  reach out if there
#     are any violations."
Set-ParaRuns $d "I can have " (
    '<w:r><w:t xml:space="preserve">I can </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">access customer contact details by slip number </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">and </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> reach out if there are any violations.</w:t></w:r>'
) | Out-Null

Write-Output "done"
